$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter labels in the same order the original author did, so the shared
# string table indices line up: 0=Intended Baudrate, 1=Actual Baudrate,
# 2=Clock Frequency (MHz)
$ws.Range("A2").Value = "Intended Baudrate"
$ws.Range("A3").Value = "Actual Baudrate"
$ws.Range("A1").Value = "Clock Frequency (MHz)"

# Numeric inputs
$ws.Range("B1").Value = 100
$ws.Range("B2").Value = 115200

# Actual baudrate formula, formatted to 4 decimal places
$ws.Range("B3").Formula = "=((B1*1000000)/B2)"
$ws.Range("B3").NumberFormat = "0.0000"

# Recomputed clock-frequency check formula
$ws.Range("B4").Formula = "=(B1 * 1000000) / B3"
$ws.Range("B4").NumberFormat = "0.00"
$ws.Range("B4").HorizontalAlignment = -4131
$ws.Range("B4").IndentLevel = 9

# Column widths to fit content
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null

$ws.Range("B4").Select()
